$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 72: date 2025-10-06 (serial 45936), station "四方坪站充电量(kw)"
$row72 = @(45936, "四方坪站充电量(kw)", 723.42100000000016, 867.16, 591.80799999999999, 427.90800000000002, 382.923, 803.00599999999997, 341.93299999999999, 174.39000000000001, 127.74000000000001, 90.864000000000004, 167.101, 57.319999999999993, 512.63499999999988, 704.94999999999993, 729.33299999999997, 360.53800000000001, 374.56400000000008, 135.208, 84.28, 159.73000000000002, 171.85300000000001, 99.47999999999999, 29.5, 75.339999999999989)
for ($c = 1; $c -le 26; $c++) { $ws.Cells.Item(72, $c).Value = $row72[$c-1] }

# New row 73: date 2025-10-06 (serial 45936), station "高岭站充电量(kw)"
$row73 = @(45936, "高岭站充电量(kw)", 354.05399999999997, 300.61099999999993, 117.774, 110.821, 0, 181.762, 172.43399999999997, 57.664000000000001, 80.317999999999998, 143.376, 49.82, 62.347999999999999, 98.288999999999987, 188.672, 173.107, 99.652000000000001, 58.768999999999998, 18.312999999999999, 89.256, 0, 47.905999999999999, 20.93, 19.039000000000001, 0)
for ($c = 1; $c -le 26; $c++) { $ws.Cells.Item(73, $c).Value = $row73[$c-1] }

# Update the view's active cell/selection to match the saved workbook state
$ws.Range("D76").Select()
